$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E1").Value = "along"

# Updated row labels (unchanged text, but included for clarity) and re-run values
$data = @(
    @{ A = "<b>All</b>";          B = 0.0345677812512016;   C = 0.00955330665890686;    D = 0.0595822558434964 },
    @{ A = "<b>Europe</b>";       B = 0.0675876450707566;   C = 0.0333189915956314;     D = 0.101856298545882 },
    @{ A = "France";              B = 0.0391366742932622;   C = -0.042578009691438;     D = 0.120851358277962 },
    @{ A = "Germany";             B = 0.0947531111645335;   C = 0.0191488053021976;     D = 0.170357417026869 },
    @{ A = "Italy";               B = 0.0175119326214762;   C = -0.0577422545335324;    D = 0.0927661197764848 },
    @{ A = "Poland";              B = 0.103853982432025;    C = -0.00231517261693412;   D = 0.210023137480985 },
    @{ A = "Spain";               B = 0.0456890966545247;   C = -0.0503919462261996;    D = 0.141770139535249 },
    @{ A = "United Kingdom";      B = 0.0834478409431803;   C = -0.000972477118733037;  D = 0.167868159005094 },
    @{ A = "Switzerland";         B = 0.0490455807248947;   C = -0.0594078178313194;    D = 0.157498979281109 },
    @{ A = "Japan";               B = 0.0101133748194984;   C = -0.0428318682783557;    D = 0.0630586179173526 },
    @{ A = "USA";                 B = 0.00777904174477409;  C = -0.0365787724650717;    D = 0.0521368559546199 }
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.A
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 4).Value = $item.D
    $ws.Cells.Item($row, 5).Value = "variant_warm_glowdonation"
    $row++
}
